$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values (B/C columns, rows 2-5)
$ws.Range("C2").Value = 12.6
$ws.Range("B3").Value = 4.0999999999999996
$ws.Range("C3").Value = 11
$ws.Range("B4").Value = 0.65
$ws.Range("C4").Value = 1.6
$ws.Range("C5").Value = 25

# Add a new (blank, but styled) row 6 - match formatting of the row above it
# so the engine reuses the same cell style and the row becomes part of the
# worksheet's used range (A1:C6)
$ws.Range("A6:C6").WrapText = $true
$ws.Range("A6:C6").VerticalAlignment = -4108

# Resize the columns to (as close as this engine's rounding allows)
# the best-fit widths, mirroring the "AutoFit" the author performed after
# the data changed (target stored widths ~21.375 / 5.125 / 5.5 chars)
$ws.Columns.Item(1).ColumnWidth = 20.714285714285715
$ws.Columns.Item(2).ColumnWidth = 4.428571428571429
$ws.Columns.Item(3).ColumnWidth = 4.714285714285714
